$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Prolong the following traces by updating their end-frame numbers
# (order matters: it determines the order new shared-string entries are appended)
$ws.Range("A132").Value = "20180419_02_040_044"
$ws.Range("A126").Value = "20180405_01_134_137"
$ws.Range("A128").Value = "20180418_01_001_004"
$ws.Range("A144").Value = "20180423_01_095_100"

# Reflect the scrolled/selected view state left by the edit
$excel.ActiveWindow.ScrollRow = 136
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A144").Select()
